$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "9.74"
$ws.Range("F2").Value = -9.26
$ws.Range("G2").Value = 28.03

# Row 3
$ws.Range("E3").Value = "4.77"
$ws.Range("F3").Value = -18.65
$ws.Range("G3").Value = 28.32

# Row 4
$ws.Range("E4").Value = "3.79"
$ws.Range("F4").Value = -26.73
$ws.Range("G4").Value = 33.21

# Row 5
$ws.Range("E5").Value = "-6.5"
$ws.Range("F5").Value = -30.32
$ws.Range("G5").Value = 17.57

# Row 10
$ws.Range("F10").Value = -0.37
$ws.Range("G10").Value = 0.28

# Row 11
$ws.Range("F11").Value = 0.63
$ws.Range("G11").Value = 1.06

# Row 13
$ws.Range("E13").Value = "-0.13"
$ws.Range("F13").Value = -0.25

# Row 14
$ws.Range("E14").Value = "0.05"
$ws.Range("F14").Value = -0.09

# Row 18
$ws.Range("E18").Value = "2.59"
$ws.Range("F18").Value = 2.1

# Row 19
$ws.Range("E19").Value = "0.48"

# Row 20
$ws.Range("E20").Value = "-0.01"
$ws.Range("F20").Value = -0.49

# Row 21
$ws.Range("E21").Value = "-0.19"
$ws.Range("F21").Value = -0.77

# Row 25
$ws.Range("E25").Value = "11.98"
$ws.Range("F25").Value = 7.78
$ws.Range("G25").Value = 16.12

# Row 26
$ws.Range("E26").Value = "-0.56"
$ws.Range("F26").Value = -5.92
$ws.Range("G26").Value = 5.07

# Row 27
$ws.Range("E27").Value = "13.82"
$ws.Range("F27").Value = 2.34
$ws.Range("G27").Value = 25.45

# Row 28
$ws.Range("E28").Value = "-4.47"
$ws.Range("F28").Value = -8.97
$ws.Range("G28").Value = -0.25

# Row 29
$ws.Range("E29").Value = "9.49"
$ws.Range("F29").Value = -8.93
$ws.Range("G29").Value = 27.64

# Row 30
$ws.Range("E30").Value = "1.95"
$ws.Range("F30").Value = -3.44
$ws.Range("G30").Value = 7.43
